# Update the cryptocurrency price/volume table (cryptos list refresh).
# Cells that hold decimal-looking text (e.g. "212.71") must stay plain text
# -- like the original inlineStr cells -- instead of being auto-coerced to
# numbers by Excel, so we briefly force a text NumberFormat for those and
# restore the default "Normal" style afterwards (no lasting style change).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.257.41"
$ws.Cells.Item(3, 4).Value = "1.592.65"
$ws.Cells.Item(3, 5).Value = "  +0.07%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "212.71"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.30%  "
$ws.Cells.Item(6, 5).Value = "  -0.37%  "
$ws.Cells.Item(7, 5).Value = "  +0.07%  "
$ws.Cells.Item(8, 5).Value = "  -0.54%  "
$ws.Cells.Item(9, 5).Value = "  -0.51%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "18.93"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.21%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.0850"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.06%  "
$ws.Cells.Item(12, 4).Value = "1.816.23"
$ws.Cells.Item(12, 5).Value = "  +0.14%  "
$ws.Cells.Item(13, 4).Value = "1.593.12"
$ws.Cells.Item(13, 5).Value = "  -0.02%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "4.01"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.11%  "
$ws.Cells.Item(15, 5).Value = "  -2.83%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "63.89"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.95%  "
$ws.Cells.Item(17, 4).Value = "26.257.15"
$ws.Cells.Item(17, 5).Value = "  -0.09%  "
$ws.Cells.Item(18, 5).Value = "  -0.74%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "214.91"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.80%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "7.38"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -1.39%  "
$ws.Cells.Item(21, 5).Value = "  +0.13%  "
$ws.Cells.Item(22, 5).Value = "  -0.15%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "9.00"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.25%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "2.10"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -2.62%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "144.65"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.00%  "
$ws.Cells.Item(27, 5).Value = "  -1.35%  "
$ws.Cells.Item(28, 5).Value = "  +0.65%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "15.12"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -1.09%  "
$ws.Cells.Item(31, 5).Value = "  -0.20%  "
$ws.Cells.Item(32, 5).Value = "  -0.55%  "
$ws.Cells.Item(33, 4).Value = "1.423.17"
$ws.Cells.Item(33, 5).Value = "  +6.20%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "2.96"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.13%  "
$ws.Cells.Item(35, 5).Value = "  -1.02%  "
$ws.Cells.Item(36, 5).Value = "  -1.27%  "
$ws.Cells.Item(37, 5).Value = "  -4.47%  "
$ws.Cells.Item(38, 5).Value = "  -0.55%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.825"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.84%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "5.77"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.23%  "
$ws.Cells.Item(42, 5).Value = "  +0.94%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.925"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -8.09%  "
$ws.Cells.Item(44, 5).Value = "  -0.33%  "
$ws.Cells.Item(45, 4).Value = "1.729.06"
$ws.Cells.Item(45, 5).Value = "  +0.27%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "60.76"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.85%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "86.64"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.35%  "
$ws.Cells.Item(48, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(48, 4).Value = "0.0₆0104"
$ws.Cells.Item(48, 5).Value = "  +0.17%  "
$ws.Cells.Item(49, 2).Value = "RenderToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "1.48"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.51%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.0500"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.74%  "
$ws.Cells.Item(51, 2).Value = "Algorand"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.0951"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -2.94%  "
